$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6341752409934998
$ws.Range("B1").Value = 1.485644817352295
$ws.Range("C1").Value = 5.921566009521484
$ws.Range("D1").Value = 2.803008317947388
$ws.Range("E1").Value = 1.86019229888916
